$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 320, pushing existing rows 320:356 down to 321:357.
$ws.Rows("320:320").Insert()

# Populate the newly inserted row 320 with a duplicate of the (now shifted)
# row 321 data, except for the Fecha (D) and Volumen (J) values.
$ws.Cells.Item(320, 1).Value = 4
$ws.Cells.Item(320, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(320, 3).Value = "Los Lagos"
$ws.Cells.Item(320, 4).Value = 45212
$ws.Cells.Item(320, 5).Value = 10
$ws.Cells.Item(320, 6).Value = 100112009
$ws.Cells.Item(320, 7).Value = "Acelga"
$ws.Cells.Item(320, 8).Value = "Sin especificar"
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 75
$ws.Cells.Item(320, 11).Value = 10000
$ws.Cells.Item(320, 12).Value = 10000
$ws.Cells.Item(320, 13).Value = 10000
$ws.Cells.Item(320, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(320, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(320, 16).Value = 833
$ws.Cells.Item(320, 17).Value = 12
$ws.Cells.Item(320, 18).Value = "Hortaliza"
